$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6606524410359556
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 6708.013860684405
$ws.Range("E2").Value = 2195978.878461985
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 2202689.208753193
